# Populate Sheet1 with the saved-query translation review table
# (original English sentence / original translation / modified translation / reason of correction),
# repeated twice (rows 2-10 and 11-19), matching the committed "open saved query" data dump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 19,4

# Header row
$data[0,0] = 'original English sentence'
$data[0,1] = 'original translation'
$data[0,2] = 'modified translation'
$data[0,3] = 'reason of correction'

# Data rows (English source, Mongolian translation; "modified translation" / "reason" columns left blank)
$data[1,0] = 'Help fashion-challenged Danny pick a swell outfit for his big date with Jamie.'
$data[1,1] = 'Загвар өмсөгч Дэннид Жэймитэй хийх том болзоондоо гоё хувцас сонгоход нь туслаарай.'
$data[1,2] = ''
$data[1,3] = ''
$data[2,0] = 'Hit ''Stop'' to select the drawer containing the costume you want Danny to wear.'
$data[2,1] = '"Зогс" дээр дарж Даннигийн өмсөхийг хүссэн хувцасны шүүгээг сонгоно уу.'
$data[2,2] = ''
$data[2,3] = ''
$data[3,0] = 'Please don''t show me this dialogue again.'
$data[3,1] = 'Энэ харилцан яриаг надад дахиж битгий үзүүлээрэй.'
$data[3,2] = ''
$data[3,3] = ''
$data[4,0] = 'Avoid'
$data[4,1] = 'Зайлсхий'
$data[4,2] = ''
$data[4,3] = ''
$data[5,0] = 'Where''s Danny?'
$data[5,1] = 'Данни хаана байна'
$data[5,2] = ''
$data[5,3] = ''
$data[6,0] = 'Find'
$data[6,1] = 'Хай'
$data[6,2] = ''
$data[6,3] = ''
$data[7,0] = 'That miserable little punk Danny is trying to hide from me again. I''ll make it worth your while if you help me find him. But I''d rather avoid letting Jamie know I''m looking for Danny so let''s avoid her.'
$data[7,1] = 'Тэр хөөрхийлөлтэй бяцхан панк Дэнни дахиад л надаас нуугдах гэж байна. Хэрэв та түүнийг олоход надад тусалбал би үүнийг үнэ цэнэтэй болгоно. Гэхдээ би Дэнниг хайж байгаагаа Жэймид мэдэгдэхээс зайлсхийсэн нь дээр, тиймээс түүнээс зайлсхийцгээе.'
$data[7,2] = ''
$data[7,3] = ''
$data[8,0] = 'Click on the spots you think Danny might be hiding. Likely spots earn you coins.'
$data[8,1] = 'Дэннигийн нуугдаж байгаа гэж бодож буй газрууд дээр дарна уу. Магадгүй цэгүүд танд зоос олох болно.'
$data[8,2] = ''
$data[8,3] = ''
$data[9,0] = 'The round ends if you find Jamie'
$data[9,1] = 'Хэрэв та Жэймиг олвол тойрог дуусна'
$data[9,2] = ''
$data[9,3] = ''
$data[10,0] = 'Help fashion-challenged Danny pick a swell outfit for his big date with Jamie.'
$data[10,1] = 'Загвар өмсөгч Дэннид Жэймитэй хийх том болзоондоо гоё хувцас сонгоход нь туслаарай.'
$data[10,2] = ''
$data[10,3] = ''
$data[11,0] = 'Hit ''Stop'' to select the drawer containing the costume you want Danny to wear.'
$data[11,1] = '"Зогс" дээр дарж Даннигийн өмсөхийг хүссэн хувцасны шүүгээг сонгоно уу.'
$data[11,2] = ''
$data[11,3] = ''
$data[12,0] = 'Please don''t show me this dialogue again.'
$data[12,1] = 'Энэ харилцан яриаг надад дахиж битгий үзүүлээрэй'
$data[12,2] = ''
$data[12,3] = ''
$data[13,0] = 'Avoid'
$data[13,1] = 'Зайлсхий'
$data[13,2] = ''
$data[13,3] = ''
$data[14,0] = 'Where''s Danny?'
$data[14,1] = 'Данни хаана байна'
$data[14,2] = ''
$data[14,3] = ''
$data[15,0] = 'Find'
$data[15,1] = 'Хай'
$data[15,2] = ''
$data[15,3] = ''
$data[16,0] = 'That miserable little punk Danny is trying to hide from me again. I''ll make it worth your while if you help me find him. But I''d rather avoid letting Jamie know I''m looking for Danny so let''s avoid her.'
$data[16,1] = 'Тэр хөөрхийлөлтэй бяцхан панк Дэнни дахиад л надаас нуугдах гэж байна. Хэрэв та түүнийг олоход надад тусалбал би үүнийг үнэ цэнэтэй болгоно. Гэхдээ би Дэнниг хайж байгаагаа Жэймид мэдэгдэхээс зайлсхийсэн нь дээр, тиймээс түүнээс зайлсхийцгээе.'
$data[16,2] = ''
$data[16,3] = ''
$data[17,0] = 'Click on the spots you think Danny might be hiding. Likely spots earn you coins.'
$data[17,1] = 'Дэннигийн нуугдаж байгаа гэж бодож буй газрууд дээр дарна уу. Магадгүй цэгүүд танд зоос олох болно.'
$data[17,2] = ''
$data[17,3] = ''
$data[18,0] = 'The round ends if you find Jamie'
$data[18,1] = 'Хэрэв та Жэймиг олвол тойрог дуусна'
$data[18,2] = ''
$data[18,3] = ''

$ws.Range("A1:D19").Value = $data

